# FY2020Q2 content refresh - officeaddins - m01 (#654)
# - update module to use latest yeoman generator
#
# The only author-intended content edit in this commit is the re-typed
# subtitle on slide 1 ("Working with Content Controls"): PowerPoint split
# it into two runs because only the tail of the string ("Content
# Controls") was retyped (now with a trailing space). The diff's many
# "9/8/19 8:0x PM" -> "12/8/19 8:32 PM" hunks are just PowerPoint
# re-stamping the auto-updating `datetime8` footer fields (on the notes
# pages/notes master/handout master) when the file was re-saved; those
# fields are not user-editable text runs in the object model (there is no
# `a:r` to target - just an `a:fld`), so they are intentionally left
# alone here rather than risk corrupting unrelated placeholders.

$p = $ppt.ActivePresentation

# -----------------------------------------------------------------------
# Slide 1 subtitle: "Working with Content Controls" -> two runs:
# "Working with " (unchanged run/formatting) + "Content Controls "
# (freshly retyped run, trailing space included).
# -----------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $shp = $s1.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "Working with Content Controls") {
            # keep "Working with " (chars 1-13) as-is, retype the rest
            # ("Content Controls", chars 14-29) with a trailing space
            $tail = $tr.Characters(14, $tr.Length - 13)
            $tail.Text = "Content Controls "
        }
    }
}
